$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 472.62964
$ws.Range("J17").Value = 472.62964
$ws.Range("L17").Value = 1417.88892
$ws.Range("N17").Value = -1753.88892

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1299.75
$ws.Range("I18").Value = 1299.75
$ws.Range("K18").Value = 1299.75
$ws.Range("M18").Value = -1015.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 603.13635
$ws.Range("I33").Value = 183.33333
$ws.Range("J33").Value = 1106.9
$ws.Range("K33").Value = 183.33333
$ws.Range("L33").Value = 1106.9
$ws.Range("M33").Value = 45.66667000000001
$ws.Range("N33").Value = -1564.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2250.75
$ws.Range("I62").Value = 1800
$ws.Range("J62").Value = 2521.2
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 2521.2
$ws.Range("M62").Value = -1176
$ws.Range("N62").Value = -3769.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2250.75
$ws.Range("I65").Value = 1800
$ws.Range("J65").Value = 2521.2
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 12606
$ws.Range("M65").Value = -5880
$ws.Range("N65").Value = -18846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 240.94737
$ws.Range("I96").Value = 245.3077
$ws.Range("K96").Value = 735.9231
$ws.Range("M96").Value = 637.0769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1686.25
$ws.Range("J97").Value = 2001.6666
$ws.Range("L97").Value = 6004.9998
$ws.Range("N97").Value = -6996.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1868.4445
$ws.Range("I100").Value = 1726.25
$ws.Range("K100").Value = 1726.25
$ws.Range("M100").Value = -1185.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 957.087
$ws.Range("I103").Value = 738.7059
$ws.Range("K103").Value = 2216.1177
$ws.Range("M103").Value = -1630.1177

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1955.625
$ws.Range("I116").Value = 1727
$ws.Range("K116").Value = 1727
$ws.Range("M116").Value = 1715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 19572.092
$ws.Range("I135").Value = 21881.809
$ws.Range("J135").Value = 4064
$ws.Range("K135").Value = 196936.281
$ws.Range("L135").Value = 36576
$ws.Range("M135").Value = -194401.281
$ws.Range("N135").Value = -41646

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2826577.2
$ws.Range("I138").Value = 1289.4
$ws.Range("J138").Value = 8774552
$ws.Range("K138").Value = 3868.2
$ws.Range("L138").Value = 26323656
$ws.Range("M138").Value = 1271.8
$ws.Range("N138").Value = -26333936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 886.91
$ws.Range("I32").Value = 737.63855
$ws.Range("J32").Value = 1615.7059
$ws.Range("K32").Value = 737.63855
$ws.Range("L32").Value = 1615.7059
$ws.Range("M32").Value = -450.63855
$ws.Range("N32").Value = -2189.7059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1511.5625
$ws.Range("I45").Value = 959.5
$ws.Range("K45").Value = 959.5
$ws.Range("M45").Value = -582.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20449760
$ws.Range("I61").Value = 23833946
$ws.Range("J61").Value = 144644.86
$ws.Range("K61").Value = 23833946
$ws.Range("L61").Value = 144644.86
$ws.Range("M61").Value = -23833734
$ws.Range("N61").Value = -145068.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2977202
$ws.Range("I97").Value = 4465296
$ws.Range("K97").Value = 4465296
$ws.Range("M97").Value = -4464800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1406.95
$ws.Range("I110").Value = 1131.5454
$ws.Range("K110").Value = 1131.5454
$ws.Range("M110").Value = 913.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 62580.06
$ws.Range("I132").Value = 44352.566
$ws.Range("J132").Value = 100692.09
$ws.Range("K132").Value = 133057.698
$ws.Range("L132").Value = 302076.27
$ws.Range("M132").Value = -130527.698
$ws.Range("N132").Value = -307136.27

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 20449760
$ws.Range("I136").Value = 23833946
$ws.Range("J136").Value = 144644.86
$ws.Range("K136").Value = 71501838
$ws.Range("L136").Value = 433934.58
$ws.Range("M136").Value = -71499288
$ws.Range("N136").Value = -439034.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 684.8461
$ws.Range("I80").Value = 309.42856
$ws.Range("J80").Value = 1122.8334
$ws.Range("K80").Value = 309.42856
$ws.Range("L80").Value = 1122.8334
$ws.Range("M80").Value = 688.5714399999999
$ws.Range("N80").Value = -3118.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 684.8461
$ws.Range("I83").Value = 309.42856
$ws.Range("J83").Value = 1122.8334
$ws.Range("K83").Value = 1547.1428
$ws.Range("L83").Value = 5614.166999999999
$ws.Range("M83").Value = 3444.8572
$ws.Range("N83").Value = -15598.167

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I99").Value = 1419.091
$ws.Range("J99").Value = 1666.6666
$ws.Range("K99").Value = 1419.091
$ws.Range("L99").Value = 1666.6666
$ws.Range("M99").Value = 78.90900000000011
$ws.Range("N99").Value = -4662.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2048
$ws.Range("I107").Value = 1986.3077
$ws.Range("J107").Value = 2137.111
$ws.Range("K107").Value = 1986.3077
$ws.Range("L107").Value = 2137.111
$ws.Range("M107").Value = -66.30770000000007
$ws.Range("N107").Value = -5977.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1974.4872
$ws.Range("I134").Value = 1272.8823
$ws.Range("J134").Value = 2516.6365
$ws.Range("K134").Value = 3818.6469
$ws.Range("L134").Value = 7549.9095
$ws.Range("M134").Value = -1283.6469
$ws.Range("N134").Value = -12619.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 32260516
$ws.Range("I58").Value = 45456916
$ws.Range("J58").Value = 2646
$ws.Range("K58").Value = 45456916
$ws.Range("L58").Value = 2646
$ws.Range("M58").Value = -45456713
$ws.Range("N58").Value = -3052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 736
$ws.Range("J107").Value = 878.25
$ws.Range("L107").Value = 878.25
$ws.Range("N107").Value = -4718.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 35358.32
$ws.Range("I132").Value = 24102.045
$ws.Range("J132").Value = 68376.734
$ws.Range("K132").Value = 72306.13499999999
$ws.Range("L132").Value = 205130.202
$ws.Range("M132").Value = -69776.13499999999
$ws.Range("N132").Value = -210190.202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 17968.672
$ws.Range("I134").Value = 1070.4375
$ws.Range("J134").Value = 68663.375
$ws.Range("K134").Value = 3211.3125
$ws.Range("L134").Value = 205990.125
$ws.Range("M134").Value = -676.3125
$ws.Range("N134").Value = -211060.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 32260516
$ws.Range("I136").Value = 45456916
$ws.Range("J136").Value = 2646
$ws.Range("K136").Value = 136370748
$ws.Range("L136").Value = 7938
$ws.Range("M136").Value = -136368198
$ws.Range("N136").Value = -13038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2473.5715
$ws.Range("I4").Value = 1900
$ws.Range("J4").Value = 2494.8147
$ws.Range("K4").Value = 5700
$ws.Range("L4").Value = 7484.4441
$ws.Range("M4").Value = -5588
$ws.Range("N4").Value = -7708.4441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 661.76666
$ws.Range("I122").Value = 258.1111
$ws.Range("J122").Value = 1267.25
$ws.Range("K122").Value = 2322.9999
$ws.Range("L122").Value = 11405.25
$ws.Range("M122").Value = 127.0000999999997
$ws.Range("N122").Value = -16305.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 960.9048
$ws.Range("I124").Value = 599
$ws.Range("K124").Value = 1797
$ws.Range("M124").Value = 3113

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1083.2
$ws.Range("J131").Value = 1261.5814
$ws.Range("L131").Value = 3784.7442
$ws.Range("N131").Value = -13864.7442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3163.5425
$ws.Range("I139").Value = 1457.3704
$ws.Range("J139").Value = 4603.125
$ws.Range("K139").Value = 4372.1112
$ws.Range("L139").Value = 13809.375
$ws.Range("M139").Value = 767.8887999999997
$ws.Range("N139").Value = -24089.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1722.5
$ws.Range("I113").Value = 1134
$ws.Range("J113").Value = 1883
$ws.Range("K113").Value = 1134
$ws.Range("L113").Value = 1883
$ws.Range("M113").Value = 1036
$ws.Range("N113").Value = -6223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 472.42856
$ws.Range("I93").Value = 475.75
$ws.Range("J93").Value = 468
$ws.Range("K93").Value = 475.75
$ws.Range("L93").Value = 468
$ws.Range("M93").Value = 772.25
$ws.Range("N93").Value = -2964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21711.857
$ws.Range("I132").Value = 10901.912
$ws.Range("J132").Value = 52520.2
$ws.Range("K132").Value = 32705.736
$ws.Range("L132").Value = 157560.6
$ws.Range("M132").Value = -30175.736
$ws.Range("N132").Value = -162620.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 642.4286
$ws.Range("I107").Value = 642
$ws.Range("J107").Value = 642.8570999999999
$ws.Range("K107").Value = 1926
$ws.Range("L107").Value = 1928.5713
$ws.Range("M107").Value = -6
$ws.Range("N107").Value = -5768.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55026.12
$ws.Range("I132").Value = 46023.137
$ws.Range("J132").Value = 85497.766
$ws.Range("K132").Value = 138069.411
$ws.Range("L132").Value = 256493.298
$ws.Range("M132").Value = -135539.411
$ws.Range("N132").Value = -261553.298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 45830.09
$ws.Range("I136").Value = 29324.258
$ws.Range("K136").Value = 87972.774
$ws.Range("M136").Value = -85422.774
